# #5: property aircraft done
#
# The "property_category" column (column I on the 建物/building sheet,
# column H on the 汽車/car sheet) was left as "land" for rows that are
# not actually land. Fix the category label on each sheet so it matches
# the sheet's own asset type.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: property_category column is I, data rows 2-3
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I3").Value = "building"

# 汽車 (car) sheet: property_category column is H, data row 2
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
